$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed entirely in the diff
$ws.Cells.Item(2, 3).ClearContents()  # C2 was 1.75539628881467
$ws.Cells.Item(2, 5).ClearContents()  # E2 was 0.2337905658324813
$ws.Cells.Item(3, 3).ClearContents()  # C3 was 2.071001150990881
$ws.Cells.Item(4, 3).ClearContents()  # C4 was 2.213911448916162

# Update cell values that changed to new recomputed precision
$ws.Cells.Item(4, 5).Value = 3.386383090739975  # E4: 3.386383090739953 -> 3.386383090739975
$ws.Cells.Item(5, 3).Value = 0.8787496612563173  # C5: 0.8787496612562951 -> 0.8787496612563173
$ws.Cells.Item(5, 5).Value = 1.013823151053028  # E5: 1.013823151053095 -> 1.013823151053028
$ws.Cells.Item(6, 3).Value = 2.533533936850585  # C6: 2.533533936850563 -> 2.533533936850585
$ws.Cells.Item(6, 5).Value = 0.9842934829757288  # E6: 0.984293482975751 -> 0.9842934829757288
$ws.Cells.Item(7, 5).Value = 4.356912452939454  # E7: 4.356912452939476 -> 4.356912452939454
$ws.Cells.Item(8, 5).Value = 3.612753212925446  # E8: 3.612753212925401 -> 3.612753212925446
$ws.Cells.Item(9, 3).Value = 2.431458940166964  # C9: 2.431458940167008 -> 2.431458940166964
$ws.Cells.Item(9, 5).Value = 5.259925231829876  # E9: 5.259925231829898 -> 5.259925231829876
$ws.Cells.Item(10, 3).Value = 1.21254482274098  # C10: 1.212544822741002 -> 1.21254482274098
$ws.Cells.Item(10, 5).Value = 2.158838189283219  # E10: 2.158838189283174 -> 2.158838189283219
$ws.Cells.Item(11, 3).Value = 1.447930496829564  # C11: 1.447930496829541 -> 1.447930496829564
$ws.Cells.Item(12, 5).Value = 1.194058515117313  # E12: 1.194058515117336 -> 1.194058515117313
$ws.Cells.Item(13, 5).Value = -0.563208905821222  # E13: -0.5632089058212553 -> -0.563208905821222
$ws.Cells.Item(14, 3).Value = 0.4712609263772816  # C14: 0.4712609263772594 -> 0.4712609263772816
$ws.Cells.Item(14, 5).Value = 1.409662779709797  # E14: 1.409662779709819 -> 1.409662779709797
$ws.Cells.Item(15, 3).Value = 0.6742451383204839  # C15: 0.6742451383205061 -> 0.6742451383204839
$ws.Cells.Item(15, 5).Value = 1.713290556413583  # E15: 1.713290556413605 -> 1.713290556413583
$ws.Cells.Item(18, 5).Value = 4.595879021798344  # E18: 4.595879021798321 -> 4.595879021798344
$ws.Cells.Item(20, 3).Value = 4.109890522944326  # C20: 4.109890522944348 -> 4.109890522944326
$ws.Cells.Item(21, 3).Value = 1.715791310593251  # C21: 1.715791310593229 -> 1.715791310593251
$ws.Cells.Item(21, 5).Value = 1.687339605296501  # E21: 1.687339605296523 -> 1.687339605296501
$ws.Cells.Item(22, 5).Value = 0.02883756256673031  # E22: 0.02883756256675252 -> 0.02883756256673031
$ws.Cells.Item(23, 3).Value = 1.862609889357336  # C23: 1.862609889357314 -> 1.862609889357336
$ws.Cells.Item(24, 5).Value = -2.079848588862143  # E24: -2.079848588862154 -> -2.079848588862143
$ws.Cells.Item(25, 5).Value = 2.147322685428366  # E25: 2.147322685428343 -> 2.147322685428366
$ws.Cells.Item(26, 5).Value = 0.9262553939923146  # E26: 0.9262553939922924 -> 0.9262553939923146
$ws.Cells.Item(28, 3).Value = 0.893498267486792  # C28: 0.8934982674867697 -> 0.893498267486792
$ws.Cells.Item(28, 5).Value = -1.194610791899986  # E28: -1.194610791899997 -> -1.194610791899986
$ws.Cells.Item(30, 5).Value = 2.928189816005689  # E30: 2.928189816005666 -> 2.928189816005689
$ws.Cells.Item(31, 3).Value = 2.306826470345391  # C31: 2.306826470345347 -> 2.306826470345391
$ws.Cells.Item(31, 5).Value = 1.40519946540949  # E31: 1.405199465409468 -> 1.40519946540949
$ws.Cells.Item(32, 5).Value = 0.8024032016000104  # E32: 0.8024032015999882 -> 0.8024032016000104
$ws.Cells.Item(33, 3).Value = 3.265677646667942  # C33: 3.265677646667919 -> 3.265677646667942
$ws.Cells.Item(33, 5).Value = 5.715169758465  # E33: 5.715169758464977 -> 5.715169758465
$ws.Cells.Item(34, 5).Value = 3.828814763561783  # E34: 3.828814763561761 -> 3.828814763561783
$ws.Cells.Item(35, 3).Value = 0.2974381310041352  # C35: 0.297438131004113 -> 0.2974381310041352
$ws.Cells.Item(35, 5).Value = -2.009776081564663  # E35: -2.009776081564674 -> -2.009776081564663
$ws.Cells.Item(36, 5).Value = 0.8023688159249032  # E36: 0.8023688159249254 -> 0.8023688159249032
$ws.Cells.Item(37, 5).Value = 6.778609849419737  # E37: 6.778609849419781 -> 6.778609849419737
$ws.Cells.Item(38, 3).Value = 2.777797690741446  # C38: 2.777797690741424 -> 2.777797690741446
$ws.Cells.Item(39, 3).Value = 2.475264839201419  # C39: 2.475264839201441 -> 2.475264839201419
$ws.Cells.Item(40, 3).Value = 0.06579575777907465  # C40: 0.06579575777909685 -> 0.06579575777907465
$ws.Cells.Item(40, 5).Value = 0.2740865344839749  # E40: 0.2740865344839527 -> 0.2740865344839749
$ws.Cells.Item(43, 3).Value = 1.076435582022328  # C43: 1.076435582022306 -> 1.076435582022328
$ws.Cells.Item(43, 5).Value = 2.51686114938241  # E43: 2.516861149382388 -> 2.51686114938241
$ws.Cells.Item(45, 3).Value = -1.650648527511434  # C45: -1.650648527511467 -> -1.650648527511434
$ws.Cells.Item(45, 5).Value = -0.2006752520846145  # E45: -0.2006752520846367 -> -0.2006752520846145
$ws.Cells.Item(46, 3).Value = -1.432689847121826  # C46: -1.432689847121871 -> -1.432689847121826
$ws.Cells.Item(46, 5).Value = 0.1752798163574321  # E46: 0.1752798163574099 -> 0.1752798163574321
$ws.Cells.Item(47, 5).Value = 1.396505962682837  # E47: 1.396505962682815 -> 1.396505962682837
$ws.Cells.Item(48, 5).Value = -1.696610696428313  # E48: -1.696610696428325 -> -1.696610696428313
$ws.Cells.Item(50, 3).Value = 2.033479419175155  # C50: 2.033479419175133 -> 2.033479419175155
$ws.Cells.Item(51, 3).Value = 3.147579643557918  # C51: 3.147579643557941 -> 3.147579643557918
$ws.Cells.Item(52, 5).Value = -1.362365718491854  # E52: -1.362365718491865 -> -1.362365718491854
$ws.Cells.Item(53, 3).Value = 2.581716327283523  # C53: 2.581716327283545 -> 2.581716327283523
